{"js": "// Rationale for Project Choice - First Draft Research Question - rewording\n//\n// 1) \"Research Question\" section: tweak wording in the two body paragraphs.\n// 2) \"Rationale for Project Choice\" section: add two new paragraphs of\n//    explanatory text right before the (empty, bookmark-only) paragraph\n//    that follows the heading.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Locate the paragraphs we need to touch, by distinctive content ---\nlet studioPara = null;          // \"In the modern day games studio...\"\nlet aimPara = null;             // \"I aim to create a simple tool...\"\nlet bookmarkPara = null;        // empty paragraph right after the\n                                 // \"Rationale for Project Choice\" heading\n\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (studioPara === null && t.indexOf(\"In the modern day games studio\") !== -1) {\n    studioPara = items[i];\n  } else if (aimPara === null && t.indexOf(\"I aim to create a simple tool\") !== -1) {\n    aimPara = items[i];\n  } else if (bookmarkPara === null && t.indexOf(\"Rationale for Project Choice\") !== -1) {\n    // The paragraph immediately after the heading is the target.\n    bookmarkPara = items[i + 1];\n  }\n}\n\nif (!studioPara || !aimPara || !bookmarkPara) {\n  throw new Error(\"Could not locate expected paragraphs in document\");\n}\n\n// --- 1) Reword the \"Research Question\" paragraphs ---\nstudioPara.insertText(\n  \"In the modern day games studio, artists and designers are often found using keyboard and mouse input to create scenes, art assets and such; for games. However, creative people have a tendency to work better with their hands. The keyboard and mouse input may limit their ability to do this. \",\n  \"Replace\"\n);\n\naimPara.insertText(\n  \"I aim to create a simple tool, where the input is based upon the user within their 3D environment as well as using other inputs such as the users\\u2019 voice. Creating an interface more in tune with its users\\u2019 tendencies. Resulting in the exploration of the users\\u2019 potential productivity gain and a potential higher quality of work.\",\n  \"Replace\"\n);\n\n// --- 2) Insert the two new \"Rationale for Project Choice\" paragraphs ---\n// Insert the second new paragraph's text first (prepended to the start of\n// the bookmark paragraph, ahead of the bookmark), then insert a whole new\n// paragraph before that one for the first paragraph of new text. Doing it\n// in this order keeps both insertions anchored on the same stable\n// `bookmarkPara` reference.\nconst startRange = bookmarkPara.getRange(\"Start\");\nstartRange.insertText(\n  \"With this, I have first-hand experience of how an artist works and how a programmer creates software, for how they think an artist works, as well as experience with user interfaces, tools graphics/rendering and the Microsoft Kinect.\",\n  \"Before\"\n);\n\nbookmarkPara.insertParagraph(\n  \"My inspiration for this project was found whilst on work placement at \\u2018Blitz Games Studios\\u2019. Whilst there I spent time working on their tool system (\\u2018Blitz Tech\\u2019) as well as working closely with game teams and at points the Microsoft Kinect. \",\n  \"Before\"\n);\n\nawait context.sync();\n", "ps1": "# Rationale for Project Choice - First Draft Research Question - rewording\n#\n# 1) \"Research Question\" section: tweak wording in the two body paragraphs.\n# 2) \"Rationale for Project Choice\" section: add two new paragraphs of\n#    explanatory text right before the (empty, bookmark-only) paragraph\n#    that follows the heading.\n\n$d = $word.ActiveDocument\n\n# --- 1) Reword the \"Research Question\" paragraphs ---\n\n# \"In the modern day games studio artists\" -> \"...games studio, artists\"\n$r = $d.Content\n$null = $r.Find.Execute(\n    \"games studio artists\", `\n    $false, $false, $false, $false, $false, $true, 1, $false, `\n    \"games studio, artists\", 2)\n\n# \"user in there 3D space\" -> \"user within their 3D environment\"\n$r = $d.Content\n$null = $r.Find.Execute(\n    \"the user in there 3D space\", `\n    $false, $false, $false, $false, $false, $true, 1, $false, `\n    \"the user within their 3D environment\", 2)\n\n# \"the users voice\" -> \"the users' voice\"\n$r = $d.Content\n$null = $r.Find.Execute(\n    \"the users voice\", `\n    $false, $false, $false, $false, $false, $true, 1, $false, `\n    \"the users\" + [char]8217 + \" voice\", 2)\n\n# \"its user's tendencies. Exploring improvements in productivity and quality of work.\"\n# -> \"its users' tendencies. Resulting in the exploration of the users' potential\n#     productivity gain and a potential higher quality of work.\"\n$r = $d.Content\n$oldTail = \"its user\" + [char]8217 + \"s tendencies. Exploring improvements in productivity and quality of work.\"\n$newTail = \"its users\" + [char]8217 + \" tendencies. Resulting in the exploration of the users\" + [char]8217 + \" potential productivity gain and a potential higher quality of work.\"\n$null = $r.Find.Execute(\n    $oldTail, `\n    $false, $false, $false, $false, $false, $true, 1, $false, `\n    $newTail, 2)\n\n# --- 2) Insert the two new \"Rationale for Project Choice\" paragraphs ---\n\n$heading = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Rationale for Project Choice\")) {\n        $heading = $p\n        break\n    }\n}\n\n# Paragraph right after the heading - currently empty, holds the _GoBack bookmark.\n$bookmarkPara = $heading.Next()\n\n# Prepend the \"With this, ...\" sentence to that paragraph, ahead of the bookmark.\n$insPt = $bookmarkPara.Range\n$insPt.Collapse(1)  # wdCollapseStart\n$withThisText = \"With this, I have first-hand experience of how an artist works and how a programmer creates software, for how they think an artist works, as well as experience with user interfaces, tools graphics/rendering and the Microsoft Kinect.\"\n$insPt.InsertBefore($withThisText)\n\n# Insert a brand new paragraph (\"My inspiration ...\") before that paragraph.\n$blitzText = \"My inspiration for this project was found whilst on work placement at \" + [char]8216 + \"Blitz Games Studios\" + [char]8217 + \". Whilst there I spent time working on their tool system (\" + [char]8216 + \"Blitz Tech\" + [char]8217 + \") as well as working closely with game teams and at points the Microsoft Kinect. \"\n$insPt2 = $bookmarkPara.Range\n$insPt2.Collapse(1)  # wdCollapseStart\n$insPt2.InsertBefore($blitzText + \"`r\")\n"}
